$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "tertiary" condition with "secondary" in the reward formulas
# (columns S and T, rows 2-301) and in the COUNTIF summary formula (S307).
$dataRange = $ws.Range("S2:T301")
$dataRange.Replace('"tertiary"', '"secondary"')

$countRange = $ws.Range("S307")
$countRange.Replace('"tertiary"', '"secondary"')

# Reflect the scrolled/selected view from the source edit (best effort).
$win = $excel.ActiveWindow
$win.ScrollRow = 274
$win.ScrollColumn = 1
$ws.Range("S308").Select()
